$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet: SCD0304 -> SCD0019
$ws.Name = "SCD0019"

# Update TC_ID value in B2: DGS-319 -> SCD0019-002
$ws.Range("B2").Value = "SCD0019-002"

# Widen column B so the longer TC_ID text fits (target stored width ~12.43 chars;
# COM ColumnWidth is quantized to the sheet's pixel grid, 11.65 is the closest
# input that lands on the nearest reachable grid value)
$ws.Columns.Item(2).ColumnWidth = 11.65

# Move the active selection from M2 to B3 (also clears the stale topLeftCell scroll position)
$ws.Range("B3").Select() | Out-Null
